# This script applies the translation-sheet edit described in the commit:
# "Changed dsiplay name of treatise"
#
# Summary of the change (from the OOXML diff of xl/worksheets/sheet1.xml):
#  1. Row 329 (errors.INVALID_DATAPACK_HASH) is renamed to
#     errors.INVALID_DATAPACK_PHYLUM, and its English text changes from
#     "...hash provided." to "...phylum provided."
#  2. Row 418 (errors.NO_MODELS) English text changes from
#     "No models available. Please add a model by clicking the chart in the
#     main view." to "At least 2 models are required for conversion. Please
#     add a model by clicking the chart in the main view."
#  3. A brand-new row is inserted right before the old row 424
#     (button.generate), re-introducing the key/value pair that row 329 used
#     to hold (errors.INVALID_DATAPACK_HASH / "Invalid treatise datapack hash
#     provided."), with no Chinese translation. All the following rows shift
#     down by one (424->425 ... 460->461), growing the sheet from 460 to 461
#     data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the existing INVALID_DATAPACK_HASH key (row 329) to
#    INVALID_DATAPACK_PHYLUM and update its English description.
$ws.Range("A329").Value = "errors.INVALID_DATAPACK_PHYLUM"
$ws.Range("B329").Value = "Invalid treatise datapack phylum provided."

# 2) Update the NO_MODELS error message (row 418).
$ws.Range("B418").Value = "At least 2 models are required for conversion. Please add a model by clicking the chart in the main view."

# 3) Insert a new row above row 424 (button.generate) and repopulate it with
#    the original INVALID_DATAPACK_HASH key/value (this pushes every
#    subsequent row down by one).
$ws.Rows.Item(424).Insert()
$ws.Range("A424").Value = "errors.INVALID_DATAPACK_HASH"
$ws.Range("B424").Value = "Invalid treatise datapack hash provided."
